$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "akhil"
$ws.Range("B2").Value = ""

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2020-09-04"
$ws.Range("B5").Style = "Normal"
